# Adds a new "2020" data column (column Q) to the worksheet, mirroring the
# formatting of the existing "2019" column (column P), and updates the
# sheet view's selection to column T (mirrors the authored workbook, which
# was presumably left with column T selected after the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column Q (2020), one per row (header row 4, data rows 5-14).
$values = @{
    4  = 2020
    5  = 4.5999999999999996
    6  = 4.2
    7  = 1.3
    8  = 10.8
    9  = 6.5
    10 = 2.9
    11 = 2.6
    12 = 13.1
    13 = 1
    14 = 1.3
}

# Copy the formatting of column P (2019) into column Q (2020) in one shot so
# each new Q cell picks up the same number format / font / borders as the
# corresponding P cell in its row, then overwrite with the 2020 values.
$ws.Range("P4:P14").Copy() | Out-Null
$ws.Range("Q4:Q14").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

foreach ($row in $values.Keys) {
    $ws.Range("Q$row").Value = $values[$row]
}

# Update the sheet view selection: column T (an empty column) is selected,
# with T1 as the active cell.
$ws.Range("T1:T1048576").Select() | Out-Null
